$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (A6) previously used the "date only" style (s=3); that style now
# belongs to the new last row (A7), so A6 reverts to the regular date/time
# style (s=2) that the rows above it use.
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New data row for 2021-11-11.
$ws.Range("A7").Value = 44511
$ws.Range("B7").Value = 43734

# A7 gets the "date only" style that A6 used to have.
$ws.Range("A7").NumberFormat = "YYYY-MM-DD"
